$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "test"
$ws.Range("B2").Value = 0.870955228805542
$ws.Range("C2").Value = 0.872588038444519
$ws.Range("D2").Value = 295.5448913574219
$ws.Range("E2").Value = 32.60151290893555
$ws.Range("F2").Value = 33.42208480834961
$ws.Range("G2").Value = 181.9390106201172
$ws.Range("H2").Value = 226.1534881591797
